$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 201.4397426666667
$ws.Range("H2").Value = 604.3192280000001
$ws.Range("I2").Value = 0.4833500233086392
$ws.Range("J2").Value = 0.4833500233086393
$ws.Range("O2").Value = 0.1516674503452159
$ws.Range("P2").Value = 0.1516674503452159
$ws.Range("Q2").Value = 163.174652029192
$ws.Range("R2").Value = 1468.571868262728
$ws.Range("S2").Value = 0.073308465659522
$ws.Range("T2").Value = 0.07330846565952201

$ws.Range("G3").Value = 201.4397426666667
$ws.Range("H3").Value = 604.3192280000001
$ws.Range("I3").Value = 0.4833500233086392
$ws.Range("J3").Value = 0.4833500233086393
$ws.Range("M3").Value = 4.530866666666666
$ws.Range("N3").Value = 13.5926
$ws.Range("O3").Value = 0.8483325496547841
$ws.Range("P3").Value = 0.8483325496547841
$ws.Range("Q3").Value = 912.6966153903112
$ws.Range("R3").Value = 8214.269538512801
$ws.Range("S3").Value = 0.4100415576491173
$ws.Range("T3").Value = 0.4100415576491173

$ws.Range("I4").Value = 0.1569674599353791
$ws.Range("J4").Value = 0.1569674599353792
$ws.Range("O4").Value = 0.1516674503452159
$ws.Range("P4").Value = 0.1516674503452159
$ws.Range("S4").Value = 0.02380685443556379
$ws.Range("T4").Value = 0.02380685443556379

$ws.Range("I5").Value = 0.1569674599353791
$ws.Range("J5").Value = 0.1569674599353792
$ws.Range("M5").Value = 4.530866666666666
$ws.Range("N5").Value = 13.5926
$ws.Range("O5").Value = 0.8483325496547841
$ws.Range("P5").Value = 0.8483325496547841
$ws.Range("Q5").Value = 296.3973569893778
$ws.Range("R5").Value = 2667.5762129044
$ws.Range("S5").Value = 0.1331606054998154
$ws.Range("T5").Value = 0.1331606054998154

$ws.Range("G6").Value = 60.43484133333334
$ws.Range("H6").Value = 181.304524
$ws.Range("I6").Value = 0.1450120099461104
$ws.Range("J6").Value = 0.1450120099461104
$ws.Range("O6").Value = 0.1516674503452159
$ws.Range("P6").Value = 0.1516674503452159
$ws.Range("Q6").Value = 48.95475974333601
$ws.Range("R6").Value = 440.5928376900241
$ws.Range("S6").Value = 0.02199360181796165
$ws.Range("T6").Value = 0.02199360181796165

$ws.Range("G7").Value = 60.43484133333334
$ws.Range("H7").Value = 181.304524
$ws.Range("I7").Value = 0.1450120099461104
$ws.Range("J7").Value = 0.1450120099461104
$ws.Range("M7").Value = 4.530866666666666
$ws.Range("N7").Value = 13.5926
$ws.Range("O7").Value = 0.8483325496547841
$ws.Range("P7").Value = 0.8483325496547841
$ws.Range("Q7").Value = 273.8222081024889
$ws.Range("R7").Value = 2464.3998729224
$ws.Range("S7").Value = 0.1230184081281487
$ws.Range("T7").Value = 0.1230184081281487

$ws.Range("G8").Value = 89.46554166666668
$ws.Range("H8").Value = 268.396625
$ws.Range("I8").Value = 0.2146705068098712
$ws.Range("J8").Value = 0.2146705068098712
$ws.Range("O8").Value = 0.1516674503452159
$ws.Range("P8").Value = 0.1516674503452159
$ws.Range("Q8").Value = 72.47084630275002
$ws.Range("R8").Value = 652.2376167247501
$ws.Range("S8").Value = 0.03255852843216848
$ws.Range("T8").Value = 0.03255852843216848

$ws.Range("G9").Value = 89.46554166666668
$ws.Range("H9").Value = 268.396625
$ws.Range("I9").Value = 0.2146705068098712
$ws.Range("J9").Value = 0.2146705068098712
$ws.Range("M9").Value = 4.530866666666666
$ws.Range("N9").Value = 13.5926
$ws.Range("O9").Value = 0.8483325496547841
$ws.Range("P9").Value = 0.8483325496547841
$ws.Range("Q9").Value = 405.3564405527778
$ws.Range("R9").Value = 3648.207964975
$ws.Range("S9").Value = 0.1821119783777027
$ws.Range("T9").Value = 0.1821119783777028
